$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1080.7916
$ws.Range("I15").Value = 1080.7916
$ws.Range("K15").Value = 3242.3748
$ws.Range("M15").Value = -3073.3748

$ws.Range("H18").Value = 1018.6667
$ws.Range("I18").Value = 933.0769
$ws.Range("J18").Value = 1575
$ws.Range("K18").Value = 933.0769
$ws.Range("L18").Value = 1575
$ws.Range("M18").Value = -649.0769
$ws.Range("N18").Value = -2143

$ws.Range("H98").Value = 856.4857
$ws.Range("I98").Value = 568.7083
$ws.Range("K98").Value = 568.7083
$ws.Range("M98").Value = 929.2917

$ws.Range("H106").Value = 4198.5
$ws.Range("I106").Value = 4310.625
$ws.Range("K106").Value = 4310.625
$ws.Range("M106").Value = -3679.625

$ws.Range("H107").Value = 276.94116
$ws.Range("I107").Value = 179.57143
$ws.Range("K107").Value = 179.57143
$ws.Range("M107").Value = 1740.42857

$ws.Range("H122").Value = 856.4857
$ws.Range("I122").Value = 568.7083
$ws.Range("K122").Value = 1706.1249
$ws.Range("M122").Value = 743.8751

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 788
$ws.Range("I2").Value = 698.8333
$ws.Range("K2").Value = 698.8333
$ws.Range("M2").Value = -585.8333

$ws.Range("H38").Value = 6668666.5
$ws.Range("I38").Value = 3000
$ws.Range("J38").Value = 20000000
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 20000000
$ws.Range("M38").Value = -2533
$ws.Range("N38").Value = -20000934

$ws.Range("H45").Value = 5436644.5
$ws.Range("I45").Value = 2074.7273
$ws.Range("J45").Value = 10418333
$ws.Range("K45").Value = 2074.7273
$ws.Range("L45").Value = 10418333
$ws.Range("M45").Value = -1697.7273
$ws.Range("N45").Value = -10419087

$ws.Range("H74").Value = 50343.81
$ws.Range("I74").Value = 84997.914
$ws.Range("K74").Value = 84997.914
$ws.Range("M74").Value = -84123.914

$ws.Range("H77").Value = 50343.81
$ws.Range("I77").Value = 84997.914
$ws.Range("K77").Value = 424989.57
$ws.Range("M77").Value = -420621.57

$ws.Range("H95").Value = 158594.33
$ws.Range("J95").Value = 158594.33
$ws.Range("L95").Value = 158594.33
$ws.Range("N95").Value = -164086.33

$ws.Range("H97").Value = 626.8
$ws.Range("I97").Value = 457.2857
$ws.Range("K97").Value = 457.2857
$ws.Range("M97").Value = 38.71429999999998

$ws.Range("H110").Value = 439.8095
$ws.Range("I110").Value = 520.9375
$ws.Range("J110").Value = 180.2
$ws.Range("K110").Value = 520.9375
$ws.Range("L110").Value = 180.2
$ws.Range("M110").Value = 1524.0625
$ws.Range("N110").Value = -4270.2

$ws.Range("H116").Value = 788
$ws.Range("I116").Value = 698.8333
$ws.Range("K116").Value = 698.8333
$ws.Range("M116").Value = 1595.1667

$ws.Range("H132").Value = 2442.4092
$ws.Range("I132").Value = 2039.5625
$ws.Range("K132").Value = 6118.6875
$ws.Range("M132").Value = -3588.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 788
$ws.Range("I3").Value = 698.8333
$ws.Range("K3").Value = 698.8333
$ws.Range("M3").Value = -584.8333

$ws.Range("H86").Value = 3972.4783
$ws.Range("I86").Value = 3177.4
$ws.Range("J86").Value = 5463.25
$ws.Range("K86").Value = 3177.4
$ws.Range("L86").Value = 5463.25
$ws.Range("M86").Value = -2054.4
$ws.Range("N86").Value = -7709.25

$ws.Range("H89").Value = 3972.4783
$ws.Range("I89").Value = 3177.4
$ws.Range("J89").Value = 5463.25
$ws.Range("K89").Value = 15887
$ws.Range("L89").Value = 27316.25
$ws.Range("M89").Value = -10271
$ws.Range("N89").Value = -38548.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1763.2931
$ws.Range("I31").Value = 1466.2565
$ws.Range("K31").Value = 1466.2565
$ws.Range("M31").Value = -1171.2565

$ws.Range("H34").Value = 1763.2931
$ws.Range("I34").Value = 1466.2565
$ws.Range("K34").Value = 1466.2565
$ws.Range("M34").Value = -1264.2565

$ws.Range("H99").Value = 4631979.5
$ws.Range("I99").Value = 5850162
$ws.Range("J99").Value = 2885.4
$ws.Range("K99").Value = 5850162
$ws.Range("L99").Value = 2885.4
$ws.Range("M99").Value = -5848664
$ws.Range("N99").Value = -5881.4

$ws.Range("H126").Value = 4631979.5
$ws.Range("I126").Value = 5850162
$ws.Range("J126").Value = 2885.4
$ws.Range("K126").Value = 17550486
$ws.Range("L126").Value = 8656.200000000001
$ws.Range("M126").Value = -17548016
$ws.Range("N126").Value = -13596.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 12.384615
$ws.Range("I12").Value = 7.25
$ws.Range("J12").Value = 14.666667
$ws.Range("K12").Value = 21.75
$ws.Range("L12").Value = 44.000001
$ws.Range("M12").Value = 151.25
$ws.Range("N12").Value = -390.000001

$ws.Range("H140").Value = 2554
$ws.Range("J140").Value = 3749.8333
$ws.Range("L140").Value = 11249.4999
$ws.Range("N140").Value = -21609.4999

$ws.Range("H141").Value = 1921.6364
$ws.Range("I141").Value = 1921.6364
$ws.Range("K141").Value = 5764.9092
$ws.Range("M141").Value = -584.9092000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I18").Value = 5999.3335
$ws.Range("J18").Value = 6998
$ws.Range("K18").Value = 5999.3335
$ws.Range("L18").Value = 6998
$ws.Range("M18").Value = -5706.3335
$ws.Range("N18").Value = -7584

$ws.Range("H21").Value = 732642.8
$ws.Range("I21").Value = 2502999.8
$ws.Range("J21").Value = 24500
$ws.Range("K21").Value = 2502999.8
$ws.Range("L21").Value = 24500
$ws.Range("M21").Value = -2502826.8
$ws.Range("N21").Value = -24846

$ws.Range("H30").Value = 732642.8
$ws.Range("I30").Value = 2502999.8
$ws.Range("J30").Value = 24500
$ws.Range("K30").Value = 2502999.8
$ws.Range("L30").Value = 24500
$ws.Range("M30").Value = -2502894.8
$ws.Range("N30").Value = -24710

$ws.Range("H104").Value = 33000
$ws.Range("J104").Value = 33000
$ws.Range("L104").Value = 33000
$ws.Range("N104").Value = -39988

$ws.Range("H114").Value = 76177.5
$ws.Range("J114").Value = 76177.5
$ws.Range("L114").Value = 76177.5
$ws.Range("N114").Value = -84855.5

$ws.Range("H122").Value = 1995.3478
$ws.Range("J122").Value = 2570.889
$ws.Range("L122").Value = 7712.667
$ws.Range("N122").Value = -12612.667

$ws.Range("H140").Value = 98662.336
$ws.Range("J140").Value = 98662.336
$ws.Range("L140").Value = 98662.336
$ws.Range("N140").Value = -109022.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5061.222
$ws.Range("I7").Value = 3540.4
$ws.Range("K7").Value = 3540.4
$ws.Range("M7").Value = -3428.4

$ws.Range("H16").Value = 848.8333
$ws.Range("I16").Value = 848.8333
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 848.8333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -678.8333
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 10327
$ws.Range("I22").Value = 1750
$ws.Range("K22").Value = 1750
$ws.Range("M22").Value = -1455

$ws.Range("H27").Value = 10327
$ws.Range("I27").Value = 1750
$ws.Range("K27").Value = 1750
$ws.Range("M27").Value = -1643

$ws.Range("H35").Value = 461
$ws.Range("I35").Value = 461
$ws.Range("K35").Value = 461
$ws.Range("M35").Value = -125

$ws.Range("H40").Value = 3090362.2
$ws.Range("I40").Value = 4192.0835
$ws.Range("K40").Value = 4192.0835
$ws.Range("M40").Value = -4056.0835

$ws.Range("H93").Value = 1285.8334
$ws.Range("I93").Value = 1285.8334
$ws.Range("K93").Value = 1285.8334
$ws.Range("M93").Value = -37.83339999999998

$ws.Range("H100").Value = 11437.315
$ws.Range("I100").Value = 14989.454
$ws.Range("J100").Value = 6553.125
$ws.Range("K100").Value = 14989.454
$ws.Range("L100").Value = 6553.125
$ws.Range("M100").Value = -14448.454
$ws.Range("N100").Value = -7635.125

$ws.Range("H126").Value = 5061.222
$ws.Range("I126").Value = 3540.4
$ws.Range("K126").Value = 10621.2
$ws.Range("M126").Value = -8151.200000000001

$ws.Range("H132").Value = 4585.5713
$ws.Range("I132").Value = 4459.8
$ws.Range("K132").Value = 13379.4
$ws.Range("M132").Value = -10849.4

$ws.Range("H136").Value = 3414.8708
$ws.Range("I136").Value = 3656.1765
$ws.Range("J136").Value = 3121.8572
$ws.Range("K136").Value = 10968.5295
$ws.Range("L136").Value = 9365.571599999999
$ws.Range("M136").Value = -8418.529500000001
$ws.Range("N136").Value = -14465.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2322.7144
$ws.Range("I132").Value = 2016.1578
$ws.Range("K132").Value = 6048.4734
$ws.Range("M132").Value = -3518.4734
